$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 489, pushing existing rows 489-610 down to 490-611.
$ws.Rows(489).Insert()

# Populate the newly inserted row 489 with the new record's data.
$ws.Range("A489").Value = 8
$ws.Range("B489").Value = "Terminal La Palmera de La Serena"
$ws.Range("C489").Value = "Coquimbo"
$ws.Range("D489").Value = 45204
$ws.Range("E489").Value = 4
$ws.Range("F489").Value = 100114013
$ws.Range("G489").Value = "Zanahoria"
$ws.Range("H489").Value = "Sin especificar"
$ws.Range("I489").Value = "Primera"
$ws.Range("J489").Value = 460
$ws.Range("K489").Value = 5800
$ws.Range("L489").Value = 6000
$ws.Range("M489").Value = 5900
$ws.Range("N489").Value = "$/saco 20 kilos"
$ws.Range("O489").Value = "Provincia del Elquí"
$ws.Range("P489").Value = 295
$ws.Range("Q489").Value = 20
$ws.Range("R489").Value = "Hortaliza"
